{"js": "// Office.js (Word JavaScript API) script.\n// Applies two edits to the CCT College Dublin assessment cover sheet:\n//   1. Fills in the (previously empty) \"Assessment Title:\" value cell with\n//      the report title.\n//   2. Updates the \"Date of Submission:\" day-of-month from 13 to 26\n//      (13/05/2023 -> 26/05/2023), leaving the surrounding \"/05/2023\" runs\n//      untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst coverTable = tables.items[0];\n\n// --- 1. Assessment Title value cell (row index 1, value column index 1) ---\nconst titleCell = coverTable.getCell(1, 1);\ntitleCell.body.insertText(\n  \"Sentiment Analysis and Distributed Data Process of a Twitter Dataset\",\n  Word.InsertLocation.replace\n);\n\n// --- 2. Date of Submission value cell (row index 6, value column index 1) ---\nconst dateCell = coverTable.getCell(6, 1);\nconst dateSearchResults = dateCell.body.search(\"13\", { matchCase: true, matchWholeWord: false });\ndateSearchResults.load(\"items\");\nawait context.sync();\n\ndateSearchResults.items[0].insertText(\"26\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies two edits to the CCT College Dublin assessment cover sheet:\n#   1. Fills in the (previously empty) \"Assessment Title:\" value cell with\n#      the report title.\n#   2. Updates the \"Date of Submission:\" day-of-month from 13 to 26\n#      (13/05/2023 -> 26/05/2023), leaving the surrounding \"/05/2023\" runs\n#      untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- 1. Assessment Title value cell (row 2, col 2 - 1-indexed COM table) ---\n$titleCell = $t.Cell(2, 2)\n$titleCell.Range.Text = \"Sentiment Analysis and Distributed Data Process of a Twitter Dataset\"\n\n# --- 2. Date of Submission value cell (row 7, col 2 - 1-indexed COM table) ---\n$dateCell = $t.Cell(7, 2)\n$dateRange = $dateCell.Range\n$find = $dateRange.Find\n$find.Text = \"13\"\n$found = $find.Execute()\nif ($found) {\n    $dateRange.Text = \"26\"\n}\n"}
